$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 23
$ws.Range("G4").Formula = "=MAX(E4:F4)*3+MIN(E4:F4)"
$ws.Range("G5:G19").Formula = "=MAX(E5:F5)*3+MIN(E5:F5)"

$ws.Range("B21").Value = 8
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 23

$ws.Range("D22").Select() | Out-Null
